# Auto-generated edit script: applies cell-level value updates to replicate the
# target diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 173.28572
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 6887.5
$ws.Range("J40").Value = 7714.2856
$ws.Range("L40").Value = 7714.2856
$ws.Range("N40").Value = -8064.2856
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H82").Value = 4481.3335
$ws.Range("I82").Value = 4481.3335
$ws.Range("K82").Value = 13444.0005
$ws.Range("M82").Value = -13038.0005
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492
$ws.Range("H85").Value = 4481.3335
$ws.Range("I85").Value = 4481.3335
$ws.Range("K85").Value = 13444.0005
$ws.Range("M85").Value = -12040.0005
$ws.Range("H130").Value = 98747.5
$ws.Range("J130").Value = 98747.5
$ws.Range("L130").Value = 98747.5
$ws.Range("N130").Value = -108787.5
$ws.Range("H137").Value = 3786.0908
$ws.Range("I137").Value = 3238.7778
$ws.Range("K137").Value = 9716.3334
$ws.Range("M137").Value = -7166.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3549.8125
$ws.Range("I74").Value = 2321.889
$ws.Range("J74").Value = 5128.5713
$ws.Range("K74").Value = 2321.889
$ws.Range("L74").Value = 5128.5713
$ws.Range("M74").Value = -1447.889
$ws.Range("N74").Value = -6876.5713
$ws.Range("H77").Value = 3549.8125
$ws.Range("I77").Value = 2321.889
$ws.Range("J77").Value = 5128.5713
$ws.Range("K77").Value = 11609.445
$ws.Range("L77").Value = 25642.8565
$ws.Range("M77").Value = -7241.445
$ws.Range("N77").Value = -34378.85649999999
$ws.Range("H88").Value = 662.5
$ws.Range("I88").Value = 662.5
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 662.5
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -256.5
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 662.5
$ws.Range("I91").Value = 662.5
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 662.5
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 741.5
$ws.Range("N91").ClearContents()
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H107").Value = 50225.5
$ws.Range("J107").Value = 50225.5
$ws.Range("L107").Value = 50225.5
$ws.Range("N107").Value = -57905.5
$ws.Range("H133").Value = 99995
$ws.Range("J133").Value = 99995
$ws.Range("L133").Value = 99995
$ws.Range("N133").Value = -105055

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 998.8333
$ws.Range("I5").Value = 993
$ws.Range("K5").Value = 993
$ws.Range("M5").Value = -880
$ws.Range("H64").Value = 695
$ws.Range("J64").Value = 695
$ws.Range("L64").Value = 695
$ws.Range("N64").Value = -1145
$ws.Range("H67").Value = 695
$ws.Range("J67").Value = 695
$ws.Range("L67").Value = 695
$ws.Range("N67").Value = -2255
$ws.Range("H86").Value = 1480
$ws.Range("J86").Value = 2250
$ws.Range("L86").Value = 2250
$ws.Range("N86").Value = -4496
$ws.Range("H89").Value = 1480
$ws.Range("J89").Value = 2250
$ws.Range("L89").Value = 11250
$ws.Range("N89").Value = -22482
$ws.Range("H112").Value = 99995
$ws.Range("J112").Value = 99995
$ws.Range("L112").Value = 99995
$ws.Range("N112").Value = -102949

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 717.5
$ws.Range("I16").Value = 623.3333
$ws.Range("K16").Value = 623.3333
$ws.Range("M16").Value = -336.3333
$ws.Range("H31").Value = 6090.8096
$ws.Range("I31").Value = 4265.7
$ws.Range("K31").Value = 4265.7
$ws.Range("M31").Value = -3970.7
$ws.Range("H34").Value = 6090.8096
$ws.Range("I34").Value = 4265.7
$ws.Range("K34").Value = 4265.7
$ws.Range("M34").Value = -4063.7
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 4900
$ws.Range("I86").Value = 4900
$ws.Range("K86").Value = 4900
$ws.Range("M86").Value = -3777
$ws.Range("H89").Value = 4900
$ws.Range("I89").Value = 4900
$ws.Range("K89").Value = 24500
$ws.Range("M89").Value = -18884
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 717.5
$ws.Range("I113").Value = 623.3333
$ws.Range("K113").Value = 623.3333
$ws.Range("M113").Value = 1546.6667
$ws.Range("H132").Value = 2250
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 875373.25
$ws.Range("I4").Value = 864043.6
$ws.Range("K4").Value = 2592130.8
$ws.Range("M4").Value = -2592018.8
$ws.Range("H97").Value = 746.8
$ws.Range("I97").Value = 746.8
$ws.Range("K97").Value = 2240.4
$ws.Range("M97").Value = -1744.4
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 16600
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H80").Value = 2488.5
$ws.Range("I80").Value = 1227.5
$ws.Range("K80").Value = 1227.5
$ws.Range("M80").Value = -229.5
$ws.Range("H83").Value = 2488.5
$ws.Range("I83").Value = 1227.5
$ws.Range("K83").Value = 6137.5
$ws.Range("M83").Value = -1145.5
$ws.Range("H107").Value = 551
$ws.Range("J107").Value = 1002
$ws.Range("L107").Value = 1002
$ws.Range("N107").Value = -4842
$ws.Range("H126").Value = 5187.1
$ws.Range("I126").Value = 5233.875
$ws.Range("K126").Value = 15701.625
$ws.Range("M126").Value = -13231.625
$ws.Range("H134").Value = 99998.336
$ws.Range("J134").Value = 99998.336
$ws.Range("L134").Value = 299995.008
$ws.Range("N134").Value = -305065.008

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2900.1765
$ws.Range("I22").Value = 2280.6
$ws.Range("J22").Value = 3158.3333
$ws.Range("K22").Value = 2280.6
$ws.Range("L22").Value = 3158.3333
$ws.Range("M22").Value = -1985.6
$ws.Range("N22").Value = -3748.3333
$ws.Range("H27").Value = 2900.1765
$ws.Range("I27").Value = 2280.6
$ws.Range("J27").Value = 3158.3333
$ws.Range("K27").Value = 2280.6
$ws.Range("L27").Value = 3158.3333
$ws.Range("M27").Value = -2173.6
$ws.Range("N27").Value = -3372.3333
$ws.Range("H29").Value = 28800
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 28800
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 28800
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -29390
$ws.Range("H46").Value = 1750.5
$ws.Range("I46").Value = 1800
$ws.Range("K46").Value = 1800
$ws.Range("M46").Value = -1612

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3249
$ws.Range("I62").Value = 3249
$ws.Range("K62").Value = 3249
$ws.Range("M62").Value = -2625
$ws.Range("H65").Value = 3249
$ws.Range("I65").Value = 3249
$ws.Range("K65").Value = 16245
$ws.Range("M65").Value = -13125
$ws.Range("H81").Value = 29079.572
$ws.Range("I81").Value = 29079.572
$ws.Range("K81").Value = 58159.144
$ws.Range("M81").Value = -57098.144
$ws.Range("H84").Value = 29079.572
$ws.Range("I84").Value = 29079.572
$ws.Range("K84").Value = 290795.72
$ws.Range("M84").Value = -285491.72

